$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (title moves from "Through 2021-12-29" to "Through 2021-12-30")
$ws.Name = "Through 2021-12-30"

# Update the row label in A14
$ws.Range("A14").Value = "December (through 12-30)"

# Row 14 - December (through 12-30) values
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 41
$ws.Range("D14").Value = 0.1087

$ws.Range("F14").Value = 88
$ws.Range("G14").Value = 0.0833

$ws.Range("I14").Value = 101
$ws.Range("J14").Value = 0.114

$ws.Range("L14").Value = 70
$ws.Range("M14").Value = 0.0667

$ws.Range("O14").Value = 59
$ws.Range("P14").Value = 0.0781

$ws.Range("R14").Value = 136
$ws.Range("S14").Value = 0.0621

$ws.Range("U14").Value = 190
$ws.Range("V14").Value = 0.0104

# Row 15 - Total values
$ws.Range("B15").Value = 38
$ws.Range("C15").Value = 299
$ws.Range("D15").Value = 0.1128

$ws.Range("F15").Value = 592
$ws.Range("G15").Value = 0.1017

$ws.Range("I15").Value = 859
$ws.Range("J15").Value = 0.0813

$ws.Range("L15").Value = 678
$ws.Range("M15").Value = 0.1044

$ws.Range("O15").Value = 539
$ws.Range("P15").Value = 0.0987

$ws.Range("R15").Value = 1336
$ws.Range("S15").Value = 0.0518

$ws.Range("U15").Value = 1732
$ws.Range("V15").Value = 0.0561
